$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "other stats" rows (player / football_player / link) appended below
# the existing table (rows 1-21 already present). The cells are written in
# the same order the sheet author typed/pasted them so the shared-string
# table comes out identical to the source edit.

# --- "Emilie" block (rows 22-25) ---
$ws.Cells.Item(22, 1).Value = "Emilie"
$ws.Cells.Item(23, 1).Value = "Emilie"
$ws.Cells.Item(24, 1).Value = "Emilie"
$ws.Cells.Item(25, 1).Value = "Emilie"

$ws.Cells.Item(22, 2).Value = "Ludovic Blas"
$ws.Cells.Item(23, 2).Value = "Pogba"
$ws.Cells.Item(24, 2).Value = "Hakimi"
$ws.Cells.Item(25, 2).Value = "Marquinhos"

$ws.Cells.Item(25, 3).Value = "https://fbref.com/en/players/d5f2f82b/Marquinhos"
$ws.Cells.Item(23, 3).Value = "https://fbref.com/en/players/867239d3/Paul-Pogba"
$ws.Cells.Item(24, 3).Value = "https://fbref.com/en/players/e42d61c7/Achraf-Hakimi"
$ws.Cells.Item(22, 3).Value = "https://fbref.com/en/players/6191093d/Ludovic-Blas"

# --- "Léna" block (rows 26-29) ---
$ws.Cells.Item(26, 1).Value = "Léna"
$ws.Cells.Item(27, 1).Value = "Léna"
$ws.Cells.Item(28, 1).Value = "Léna"
$ws.Cells.Item(29, 1).Value = "Léna"

$ws.Cells.Item(28, 2).Value = "Maupay"
$ws.Cells.Item(29, 3).Value = "https://fbref.com/en/players/76a82373/Carlos-Gomez"
$ws.Cells.Item(29, 2).Value = "Carlos Gomez"
$ws.Cells.Item(28, 3).Value = "https://fbref.com/en/players/4bcf39f6/Neal-Maupay"

$ws.Cells.Item(26, 2).Value = "Tolisso"
$ws.Cells.Item(27, 2).Value = "Vitinha"
$ws.Cells.Item(26, 3).Value = "https://fbref.com/en/players/652d4c37/Corentin-Tolisso"
$ws.Cells.Item(27, 3).Value = "https://fbref.com/en/players/3b029691/Vitinha"

# Widen column B so the new football-player names fit (target stored width
# 17.83203125 character-units; this engine quantizes ColumnWidth to whole
# pixel steps, so 16.95 lands on the closest reachable bucket, 17.8333...).
$ws.Columns.Item(2).ColumnWidth = 16.95

# The "Maupay" link cell (C28) got an explicit black font colour applied.
$ws.Range("C28").Font.Color = 0

# Leave the selection where the author left it after typing the new rows.
$ws.Range("A29").Select() | Out-Null
